$d = $word.ActiveDocument

# --- Edit 1: "Arial underline 16 size" paragraph -------------------------
# Add <w:lang w:val="en-US"/> to the run holding the "6" (it currently has
# none), while keeping its own run attributes (w:rsidR="007D7A0C") intact.
# We do this via InsertXML (while the old _GoBack bookmark is still present,
# so the freshly inserted run does not get silently merged with its
# neighbours), and only afterwards remove the stale _GoBack bookmark that
# used to sit right after this run.

$rFind = $d.Content
$rFind.Find.Execute("Arial underline 1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sixPos = $rFind.End

$sixRange = $d.Range($sixPos, $sixPos + 1)
$sixXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r w:rsidR="007D7A0C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>6</w:t></w:r></w:p>'
$sixRange.InsertXML($sixXml)

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# --- Edit 2: "Heading style calibri light 15 blue" -> "... 16 blue" ------
# Plain text replacement 15 -> 16 (keeps formatting/run intact).
$d.Content.Find.Execute("light 15 blue", $true, $false, $false, $false, $false, $true, 1, $false, "light 16 blue", 2) | Out-Null

# Split " light 16 blue" so "6" becomes its own run, matching a real edit
# where the user retyped just the "6" character.
$rFind2 = $d.Content
$rFind2.Find.Execute("light 1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pos = $rFind2.End

$sixRange2 = $d.Range($pos, $pos + 1)
$sixRange2.Font.Bold = 1
$sixRange2.Font.Bold = 0

# Move the _GoBack bookmark to right after the new "6" (between "6" and " blue"),
# reflecting that this is where the user's last edit occurred.
$bmRange = $d.Range($pos + 1, $pos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
